$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, F, G, K, M across rows 2-25
# (columns A, E, H, I, J, L, N, O are unchanged)
$data = @{
    2 = @{ "B"="0.3341098981171911"; "C"="0.008798438739674452"; "D"="0.206802659972027"; "F"="3.672583135881155"; "G"="0.002523304540212878"; "K"="0.2653096789746172"; "M"="0.8704095663614595" }
    3 = @{ "B"="0.3199864335167604"; "C"="0.0077503197459734"; "D"="0.1943692210197128"; "F"="3.447555269781219"; "G"="0.002528109740544025"; "K"="0.2526750490284684"; "M"="0.7758949721506099" }
    4 = @{ "B"="0.3117317880673056"; "C"="0.007134244760571562"; "D"="0.1866733306180492"; "F"="3.309488888319578"; "G"="0.002531206939190067"; "K"="0.2453169343437764"; "M"="0.7184347792257029" }
    5 = @{ "B"="0.3084724372897227"; "C"="0.006889948497232012"; "D"="0.1835211339449216"; "F"="3.253245242944899"; "G"="0.002532506130452861"; "K"="0.2424183630181886"; "M"="0.6951560736520719" }
    6 = @{ "B"="0.3079375228522849"; "C"="0.00684978719024798"; "D"="0.1829967283893978"; "F"="3.243907032924312"; "G"="0.002532724102882204"; "K"="0.2419430761643042"; "M"="0.6912987197536609" }
    7 = @{ "B"="0.3116874088258044"; "C"="0.007130922923778371"; "D"="0.1866308847277622"; "F"="3.30873029571859"; "G"="0.00253122431024782"; "K"="0.2452774393569683"; "M"="0.7181202898653964" }
    8 = @{ "B"="0.329153266463976"; "C"="0.008431254439749125"; "D"="0.2025280882402001"; "F"="3.594968719755599"; "G"="0.002524930996756054"; "K"="0.2608699748047911"; "M"="0.8376978531240837" }
    9 = @{ "B"="0.3667366349352505"; "C"="0.01120610493582319"; "D"="0.2332384418240849"; "F"="4.157388103273945"; "G"="0.002513747755169466"; "K"="0.2946456681384859"; "M"="1.077065027695227" }
    10 = @{ "B"="0.3964164753866442"; "C"="0.01339170904117282"; "D"="0.2555556427807346"; "F"="4.571728866033169"; "G"="0.002506227811439577"; "K"="0.3214537127178971"; "M"="1.256390959914938" }
    11 = @{ "B"="0.410375639598584"; "C"="0.01442005642113031"; "D"="0.2656631323887098"; "F"="4.760571031368386"; "G"="0.002502955958681321"; "K"="0.3340918723357333"; "M"="1.338833229367737" }
    12 = @{ "B"="0.415728023498076"; "C"="0.01481454354794209"; "D"="0.2694847774155278"; "F"="4.832139512448748"; "G"="0.002501738261996535"; "K"="0.338942064193219"; "M"="1.3701854701104" }
    13 = @{ "B"="0.4145723320933712"; "C"="0.01472935527505825"; "D"="0.2686619681704485"; "F"="4.816723268261967"; "G"="0.00250199957044577"; "K"="0.3378946135239858"; "M"="1.363427142806671" }
    14 = @{ "B"="0.4108146505126626"; "C"="0.0144524084363411"; "D"="0.2659776551130335"; "F"="4.766457810618704"; "G"="0.002502855352218088"; "K"="0.3344896059832081"; "M"="1.341409878491334" }
    15 = @{ "B"="0.4085216192546284"; "C"="0.01428343640915131"; "D"="0.2643326926343263"; "F"="4.735676530021124"; "G"="0.002503382310748715"; "K"="0.3324123472999645"; "M"="1.327941267278746" }
    16 = @{ "B"="0.3955134662990076"; "C"="0.01332520488519862"; "D"="0.2548942448962066"; "F"="4.559395337404794"; "G"="0.002506444620758596"; "K"="0.3206367527415352"; "M"="1.251021320558124" }
    17 = @{ "B"="0.38765098280507"; "C"="0.01274621627793238"; "D"="0.2490930265124405"; "F"="4.451348631875021"; "G"="0.002508361310710924"; "K"="0.3135267645714492"; "M"="1.204061279909297" }
    18 = @{ "B"="0.3831717337671137"; "C"="0.01241639717661513"; "D"="0.2457520746107207"; "F"="4.389236164430088"; "G"="0.00250947777376731"; "K"="0.3094789470633117"; "M"="1.177132089164914" }
    19 = @{ "B"="0.3816625135201264"; "C"="0.01230527052638308"; "D"="0.2446201388788012"; "F"="4.368211437005186"; "G"="0.002509858203720218"; "K"="0.3081155600406476"; "M"="1.16802797832095" }
    20 = @{ "B"="0.3884834988836872"; "C"="0.01280751820944914"; "D"="0.2497110111617644"; "F"="4.462846909461376"; "G"="0.002508155824274466"; "K"="0.3142793184185422"; "M"="1.209051813403448" }
    21 = @{ "B"="0.4119165675392367"; "C"="0.01453361535508435"; "D"="0.2667662567563411"; "F"="4.781220363769364"; "G"="0.002502603411992059"; "K"="0.3354879860316657"; "M"="1.347873199421045" }
    22 = @{ "B"="0.4276184405725587"; "C"="0.01569137701471135"; "D"="0.2778790719863764"; "F"="4.989637426376305"; "G"="0.002499098575609603"; "K"="0.3497247494021281"; "M"="1.439381231670254" }
    23 = @{ "B"="0.4192024615998378"; "C"="0.01507068740270512"; "D"="0.2719508476398005"; "F"="4.878367785786963"; "G"="0.002500957875640861"; "K"="0.3420917161936643"; "M"="1.39046737120762" }
    24 = @{ "B"="0.3881069907072572"; "C"="0.01277979412098773"; "D"="0.2494316382694421"; "F"="4.457648523031168"; "G"="0.002508248679367469"; "K"="0.3139389648782469"; "M"="1.206795379680372" }
    25 = @{ "B"="0.3562089081100339"; "C"="0.01043034996190784"; "D"="0.2249763288466369"; "F"="4.005073464642038"; "G"="0.002516650138229248"; "K"="0.2851616062607434"; "M"="1.011737943756032" }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = [double]$rowData[$col]
    }
}
